$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in previously-empty "Week 2" (G) values for existing rows 3-11 ---
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0

# Row 7 now has a Week 1 value and a Week 2 value
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0

# Row 8 gets an Actual Time value plus Week 1 / Week 2 values
$ws.Range("D8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0

# Row 9 gets an Actual Time value plus Week 1 / Week 2 values
$ws.Range("D9").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0

# Row 10 gets an Actual Time value plus Week 1 / Week 2 values
$ws.Range("D10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0

# Row 11 just needs the Week 2 value
$ws.Range("G11").Value = 0

# --- New backlog items added for Sprint 1 feedback fixes (rows 12-14) ---
$ws.Range("B12").Value = "Fix issues regarding testing from Sprint 1 feedback"
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = "Destiny"
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 0

$ws.Range("B13").Value = "Fix issues within web application and database from Sprint 1 feedback"
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = "Matthew"
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 0

$ws.Range("B14").Value = "Fix issues within desktop application and documentation from Sprint 1 feedback"
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = "Janera"
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 0

# --- Update the current selection / scroll position to match the saved view ---
[void]$ws.Range("B16").Select()
